# Fruta / hortaliza, semanal
# Weekly refresh of the Piña (Terminal Hortofrutícola Agro Chillán) price
# series: two new observations (dated 2023-07-25) are inserted at the top
# of the data block (rows 324-325), pushing every existing record down by
# two rows (old row 324 -> new row 326, ... old row 351 -> new row 353).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows right before the first data row of the block
# (row 324), shifting the existing rows (324-351) down to (326-353).
$ws.Rows.Item(324).Insert()
$ws.Rows.Item(324).Insert()

# New row for "Primera" quality.
$newRow324 = @(
    7,
    "Terminal Hortofrutícola Agro Chillán",
    "Ñuble",
    "2023-07-25",
    16,
    "Fruta",
    100108,
    "Tropicales y subtropicales",
    100108005,
    "Piña",
    "Caramelo",
    "Primera",
    50,
    18000,
    18000,
    18000,
    "`$/caja 12 unidades",
    "Ecuador",
    1500,
    12
)

# New row for "Segunda" quality.
$newRow325 = @(
    7,
    "Terminal Hortofrutícola Agro Chillán",
    "Ñuble",
    "2023-07-25",
    16,
    "Fruta",
    100108,
    "Tropicales y subtropicales",
    100108005,
    "Piña",
    "Caramelo",
    "Segunda",
    30,
    18000,
    18000,
    18000,
    "`$/caja 14 unidades",
    "Ecuador",
    1286,
    14
)

for ($i = 0; $i -lt $newRow324.Length; $i++) {
    $ws.Cells.Item(324, $i + 1).Value() = $newRow324[$i]
}

for ($i = 0; $i -lt $newRow325.Length; $i++) {
    $ws.Cells.Item(325, $i + 1).Value() = $newRow325[$i]
}
